$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in shared string used by Module column of row 192:
# "Stock & sale rewports" -> "Stock & sale reports"
$ws.Range("C192").Value = "Stock & sale reports"

# Rows whose Status moves from "Pending" to "Done".
# These rows (33 and 38) are currently shown (not part of the already
# highlighted/yellow "new" block), so give them the same yellow row
# highlight used by the rest of the "Done" rows in that block before
# marking them hidden by the Status filter.
$newlyHighlighted = @(33, 38)
foreach ($r in $newlyHighlighted) {
    $ws.Range("A" + $r + ":F" + $r).Interior.Color = 65535
    $ws.Range("E" + $r).Value = "Done"
    $ws.Rows.Item($r).Hidden = $true
}

# These rows already have the yellow row highlight; only their Status
# flips to "Done" and they become hidden by the autofilter.
$alreadyHighlighted = @(196, 197, 213, 214, 223, 224, 225)
foreach ($r in $alreadyHighlighted) {
    $ws.Range("E" + $r).Value = "Done"
    $ws.Rows.Item($r).Hidden = $true
}

# Update the remembered selection to match the author's final cursor position.
$ws.Range("D189").Select()
